$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header text updates (Volume/Number and report date range)
$ws.Range("A8").Value = "Volume 32   Number  17"
$ws.Range("C9").Value = "Report Covering the Week  4/21/2025  Through  4/27/2025"

# Weekly crime-complaint statistics refresh
$ws.Range("C15").Value = "0"
$ws.Range("D15").Value = "0"
$ws.Range("E15").Value = "***.*"
$ws.Range("F15").Value = 1
$ws.Range("H15").Value = -50
$ws.Range("N15").Value = -66.666666666666
$ws.Range("C16").Value = 2
$ws.Range("D16").Value = 4
$ws.Range("E16").Value = -50
$ws.Range("G16").Value = 15
$ws.Range("H16").Value = -66.666666666666
$ws.Range("I16").Value = 34
$ws.Range("J16").Value = 43
$ws.Range("K16").Value = -20.930232558139
$ws.Range("L16").Value = -22.727272727272
$ws.Range("M16").Value = -61.363636363636
$ws.Range("N16").Value = -88.513513513513
$ws.Range("C17").Value = 3
$ws.Range("D17").Value = 12
$ws.Range("E17").Value = -75
$ws.Range("F17").Value = 11
$ws.Range("G17").Value = 23
$ws.Range("H17").Value = -52.173913043478
$ws.Range("I17").Value = 70
$ws.Range("J17").Value = 94
$ws.Range("K17").Value = -25.531914893617
$ws.Range("L17").Value = -4.109589041095
$ws.Range("M17").Value = 150
$ws.Range("N17").Value = 1.449275362318
$ws.Range("C18").Value = 4
$ws.Range("D18").Value = 11
$ws.Range("E18").Value = -63.636363636363
$ws.Range("F18").Value = 18
$ws.Range("G18").Value = 22
$ws.Range("H18").Value = -18.181818181818
$ws.Range("I18").Value = 81
$ws.Range("J18").Value = 77
$ws.Range("K18").Value = 5.194805194805
$ws.Range("L18").Value = -7.954545454545
$ws.Range("M18").Value = 6.578947368421
$ws.Range("N18").Value = -83.503054989816
$ws.Range("C19").Value = 10
$ws.Range("D19").Value = 7
$ws.Range("E19").Value = 42.857142857142
$ws.Range("F19").Value = 43
$ws.Range("G19").Value = 44
$ws.Range("H19").Value = -2.272727272727
$ws.Range("I19").Value = 148
$ws.Range("J19").Value = 191
$ws.Range("K19").Value = -22.513089005235
$ws.Range("L19").Value = -24.489795918367
$ws.Range("M19").Value = -1.986754966887
$ws.Range("N19").Value = -20.855614973262
$ws.Range("C20").Value = 9
$ws.Range("D20").Value = 10
$ws.Range("E20").Value = -10
$ws.Range("F20").Value = 28
$ws.Range("G20").Value = 24
$ws.Range("H20").Value = 16.666666666666
$ws.Range("I20").Value = 112
$ws.Range("J20").Value = 97
$ws.Range("K20").Value = 15.463917525773
$ws.Range("L20").Value = 34.939759036144
$ws.Range("M20").Value = 80.645161290322
$ws.Range("N20").Value = -92.802056555269
$ws.Range("C21").Value = 28
$ws.Range("D21").Value = 44
$ws.Range("E21").Value = -36.363636363636
$ws.Range("F21").Value = 106
$ws.Range("G21").Value = 131
$ws.Range("H21").Value = -19.083969465648
$ws.Range("I21").Value = 452
$ws.Range("J21").Value = 513
$ws.Range("K21").Value = -11.890838206627
$ws.Range("L21").Value = -7.942973523421
$ws.Range("M21").Value = 10.78431372549
$ws.Range("N21").Value = -82.701875239188
$ws.Range("C22").Value = 1
$ws.Range("E22").Value = -50
$ws.Range("F22").Value = 1
$ws.Range("G22").Value = 4
$ws.Range("H22").Value = -75
$ws.Range("I22").Value = 6
$ws.Range("J22").Value = 17
$ws.Range("K22").Value = -64.705882352941
$ws.Range("L22").Value = -45.454545454545
$ws.Range("M22").Value = -25
$ws.Range("C23").Value = 2
$ws.Range("D23").Value = 2
$ws.Range("E23").Value = 0
$ws.Range("F23").Value = 3
$ws.Range("H23").Value = 50
$ws.Range("I23").Value = 21
$ws.Range("J23").Value = 14
$ws.Range("K23").Value = 50
$ws.Range("L23").Value = -25
$ws.Range("M23").Value = 110
$ws.Range("C24").Value = 22
$ws.Range("D24").Value = 20
$ws.Range("E24").Value = 10
$ws.Range("F24").Value = 99
$ws.Range("G24").Value = 86
$ws.Range("H24").Value = 15.116279069767
$ws.Range("I24").Value = 343
$ws.Range("J24").Value = 380
$ws.Range("K24").Value = -9.736842105263
$ws.Range("L24").Value = -29.568788501026
$ws.Range("M24").Value = 29.924242424242
$ws.Range("C25").Value = 10
$ws.Range("E25").Value = 42.857142857142
$ws.Range("F25").Value = 36
$ws.Range("G25").Value = 35
$ws.Range("H25").Value = 2.857142857142
$ws.Range("I25").Value = 124
$ws.Range("J25").Value = 171
$ws.Range("K25").Value = -27.485380116959
$ws.Range("L25").Value = -28.323699421965
$ws.Range("C26").Value = 10
$ws.Range("D26").Value = 7
$ws.Range("E26").Value = 42.857142857142
$ws.Range("F26").Value = 36
$ws.Range("G26").Value = 44
$ws.Range("H26").Value = -18.181818181818
$ws.Range("I26").Value = 175
$ws.Range("J26").Value = 144
$ws.Range("K26").Value = 21.527777777777
$ws.Range("L26").Value = 17.44966442953
$ws.Range("M26").Value = 35.658914728682
$ws.Range("C27").Value = "0"
$ws.Range("D27").Value = "0"
$ws.Range("E27").Value = "***.*"
$ws.Range("F27").Value = 1
$ws.Range("H27").Value = -66.666666666666
$ws.Range("C28").Value = 2
$ws.Range("D28").Value = 1
$ws.Range("E28").Value = 100
$ws.Range("F28").Value = 7
$ws.Range("G28").Value = 4
$ws.Range("H28").Value = 75
$ws.Range("I28").Value = 21
$ws.Range("J28").Value = 14
$ws.Range("K28").Value = 50
$ws.Range("L28").Value = 16.666666666666
$ws.Range("L31").Value = -62.5
